$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestData")

# Update the Nacubo GL Account Category values (was " J02 - J02"), now full description.
$newAgencyClass = " J02 - STATE HIGHWAY ADMINISTRATION                       "
$ws.Range("D2").Value = $newAgencyClass
$ws.Range("D3").Value = $newAgencyClass

# Move active selection to D3, matching the saved view state.
$ws.Range("D3").Select()
